$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# New row appended below the existing header(1)/data(2) rows.
$ws.Cells.Item(3, 1).Value = 108999614
$ws.Cells.Item(3, 2).Value = 56404
$ws.Cells.Item(3, 3).Value = "Ovaliderad"
$ws.Cells.Item(3, 4).Value = "NT"
$ws.Cells.Item(3, 5).Value = 100048
$ws.Cells.Item(3, 6).Value = "Mindre hackspett"
$ws.Cells.Item(3, 7).Value = "Dryobates minor"
$ws.Cells.Item(3, 8).Value = "(Linnaeus, 1758)"

# Column I ("Antal") is stored as text "2", not a number - leading
# apostrophe forces Excel to keep it as text.
$ws.Cells.Item(3, 9).Value = "'2"

$ws.Cells.Item(3, 11).Value = "adult"
$ws.Cells.Item(3, 13).Value = "spel/sång"
$ws.Cells.Item(3, 16).Value = "Björnåsen, Hl"
$ws.Cells.Item(3, 17).Value = 357279
$ws.Cells.Item(3, 18).Value = 6339160
$ws.Cells.Item(3, 19).Value = 25
$ws.Cells.Item(3, 20).Value = "Halland"
$ws.Cells.Item(3, 21).Value = "Varberg"
$ws.Cells.Item(3, 22).Value = "Halland"
$ws.Cells.Item(3, 23).Value = "Nösslinge"

# Start/end date columns hold plain text dates, not Excel date serials.
$ws.Cells.Item(3, 25).Value = "'2023-05-12"
$ws.Cells.Item(3, 26).Value = "09:12"
$ws.Cells.Item(3, 27).Value = "'2023-05-12"
$ws.Cells.Item(3, 28).Value = "09:12"

$ws.Cells.Item(3, 30).Value = $false
$ws.Cells.Item(3, 31).Value = $false
$ws.Cells.Item(3, 33).Value = $false

# AT3 ("Bestämningsår") and AY3 ("Projektnamn") are present but empty
# text cells (not blank cells) - a lone apostrophe yields empty text.
$ws.Cells.Item(3, 46).Value = "'"
$ws.Cells.Item(3, 49).Value = "Maria Torbjörnsson"
$ws.Cells.Item(3, 50).Value = "Maria Torbjörnsson"
$ws.Cells.Item(3, 51).Value = "'"
